$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 'P. point'
$ws.Range("C8").Value = 85
$ws.Range("D8").Value = '''4'
$ws.Range("E8").Value = 'Long point  (up to 10 mtr.)'
$ws.Range("F8").Value = 662
$ws.Range("G8").Value = '''56270.00'
$ws.Range("C9").Value = 33
$ws.Range("D9").Value = '''6'
$ws.Range("E9").Value = 'On board'
$ws.Range("F9").Value = 136
$ws.Range("G9").Value = '''4488.00'
$ws.Range("A10").Value = 'Each'
$ws.Range("C10").Value = 60
$ws.Range("D10").Value = '''4.0'
$ws.Range("E10").Value = 'P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F10").Value = 50
$ws.Range("G10").Value = '''3000.00'
$ws.Range("C11").Value = 70
$ws.Range("D11").Value = '''6.0'
$ws.Range("E11").Value = 'Providing & Fixing of  3/6 pin 16 amp flush type non modular socket  made out from Industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F11").Value = 78
$ws.Range("G11").Value = '''5460.00'
$ws.Range("A12").Value = 'R. mtr.'
$ws.Range("C12").Value = 44
$ws.Range("D12").Value = '''16'
$ws.Range("E12").Value = '20 mm'
$ws.Range("F12").Value = 40
$ws.Range("G12").Value = '''1760.00'
$ws.Range("C13").Value = 55
$ws.Range("G13").Value = '''3080.00'
$ws.Range("C14").Value = 34
$ws.Range("D14").Value = '''14.0'
$ws.Range("E14").Value = 'Supply & Laying following size earth wire in horizontal or vertical run in ground/surface/recess including riveting, soldering, saddles,  making connection with GI/Cu purity purity >95%  thimble etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range("A15").Value = 'Mtr.'
$ws.Range("C15").Value = 52
$ws.Range("D15").Value = '''23'
$ws.Range("E15").Value = '8 SWG G.I. ( Hot Dipped  ) Wire '
$ws.Range("F15").Value = 20
$ws.Range("G15").Value = '''1040.00'
$ws.Range("A16").Value = 'Each'
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 66
$ws.Range("D16").Value = '''32'
$ws.Range("E16").Value = ' 50/63 A rating'
$ws.Range("F16").Value = 900
$ws.Range("G16").Value = '''59400.00'
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = ''''
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 67
$ws.Range("D17").Value = '''34'
$ws.Range("E17").Value = 'Metal door (single phase) IK-09 and IP-43 with Metal end box'
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = '''0.00'
$ws.Range("H17").Value = 0
$ws.Range("E19").Value = 'Grand Total Rs.'
$ws.Range("G19").Value = '''134498.00'
$ws.Range("H19").Value = '''134498.00'
$ws.Range("A20").Value = ''''
$ws.Range("B20").Value = ''''
$ws.Range("C20").Value = ''''
$ws.Range("D20").Value = ''''
$ws.Range("E20").Value = 'Tender Premium @ 0%'
$ws.Range("F20").Value = ''''
$ws.Range("G20").Value = '''0.00'
$ws.Range("H20").Value = '''0.00'
$ws.Range("I20").Value = ''''
$ws.Range("A21").Value = ''''
$ws.Range("B21").Value = ''''
$ws.Range("C21").Value = ''''
$ws.Range("D21").Value = ''''
$ws.Range("E21").Value = 'NET PAYABLE AMOUNT Rs.'
$ws.Range("F21").Value = ''''
$ws.Range("G21").Value = '''134498.00'
$ws.Range("H21").Value = '''134498.00'
$ws.Range("I21").Value = ''''

$ws.Range("B18").Value = ""
$ws.Range("C18").Value = ""
$ws.Range("D18").Value = ""
$ws.Range("E18").Value = ""
$ws.Range("F18").Value = ""
$ws.Range("G18").Value = ""
$ws.Range("H18").Value = ""
$ws.Range("I18").Value = ""
